$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the worker identity (prolificid/name/gender) between row 7 and row 8
$ws.Range("B7").Value = 16
$ws.Range("C7").Value = "60863a15760523386e761cfb"
$ws.Range("D7").Value = "Roshni"
$ws.Range("E7").Value = "female"

$ws.Range("B8").Value = 13
$ws.Range("C8").Value = "5697d4ae7183b8000d0fc201"
$ws.Range("D8").Value = "Tu"
$ws.Range("E8").Value = "male"

# Swap the worker identity (prolificid/name/gender) between row 9 and row 10
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "5c27de12a2b00a00018b2c16"
$ws.Range("D9").Value = "Ankai"
$ws.Range("E9").Value = "male"

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "60bd88b8fc436774352f53b9"
$ws.Range("D10").Value = "Annes"
$ws.Range("E10").Value = "female"

# Recomputed realeffort values for every data row
$ws.Range("F2").Value = 11.45740717551576
$ws.Range("F3").Value = 10.08421220545006
$ws.Range("F4").Value = 8.474440037169567
$ws.Range("F5").Value = 8.469079889133782
$ws.Range("F6").Value = 7.167398355129854
$ws.Range("F7").Value = 6.194924391488313
$ws.Range("F8").Value = 6.053014121671316
$ws.Range("F9").Value = 5.405704693608066
$ws.Range("F10").Value = 5.190617474046819
$ws.Range("F11").Value = 4.0828279198773
$ws.Range("F12").Value = 2.154970194355574
$ws.Range("F13").Value = 0.2602112850569155
